$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "308.18"
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = "1.86%"
$cell = $ws.Cells.Item(2, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "36.31"
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = "1.93%"
$cell = $ws.Cells.Item(3, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.073"
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = "1.19%"
$cell = $ws.Cells.Item(4, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.08163"
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = "3.73%"
$cell = $ws.Cells.Item(5, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.958"
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = "6.26%"
$cell = $ws.Cells.Item(6, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(7, 2)
$cell.NumberFormat = "@"
$cell.Value = "KuCoinToken"
$cell = $ws.Cells.Item(7, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.891"
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = "1.04%"
$cell = $ws.Cells.Item(7, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(8, 2)
$cell.NumberFormat = "@"
$cell.Value = "MXToken"
$cell = $ws.Cells.Item(8, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.9276"
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = "0.32%"
$cell = $ws.Cells.Item(8, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(9, 2)
$cell.NumberFormat = "@"
$cell.Value = "LiechtensteinCryptoassetsExchange"
$cell = $ws.Cells.Item(9, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.1430"
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = "8.57%"
$cell = $ws.Cells.Item(9, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(10, 2)
$cell.NumberFormat = "@"
$cell.Value = "WazirX"
$cell = $ws.Cells.Item(10, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.1946"
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = "3.29%"
$cell = $ws.Cells.Item(10, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(11, 2)
$cell.NumberFormat = "@"
$cell.Value = "MandalaExchangeToken"
$cell = $ws.Cells.Item(11, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.09256"
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = "1.47%"
$cell = $ws.Cells.Item(11, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(12, 2)
$cell.NumberFormat = "@"
$cell.Value = "BitrueCoin"
$cell = $ws.Cells.Item(12, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.03508"
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = "-0.28%"
$cell = $ws.Cells.Item(12, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(13, 2)
$cell.NumberFormat = "@"
$cell.Value = "BitMartToken"
$cell = $ws.Cells.Item(13, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.09882"
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = "0.54%"
$cell = $ws.Cells.Item(13, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(14, 2)
$cell.NumberFormat = "@"
$cell.Value = "BitForexToken"
$cell = $ws.Cells.Item(14, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.001405"
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = "-1.39%"
$cell = $ws.Cells.Item(14, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(15, 2)
$cell.NumberFormat = "@"
$cell.Value = "TigerCash"
$cell = $ws.Cells.Item(15, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.006188"
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = "-0.12%"
$cell = $ws.Cells.Item(15, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(16, 2)
$cell.NumberFormat = "@"
$cell.Value = "LEO"
$cell = $ws.Cells.Item(16, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.922"
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = "7.26%"
$cell = $ws.Cells.Item(16, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(17, 2)
$cell.NumberFormat = "@"
$cell.Value = "GateToken"
$cell = $ws.Cells.Item(17, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.171"
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = "0.45%"
$cell = $ws.Cells.Item(17, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.484"
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = "7.42%"
$cell = $ws.Cells.Item(18, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = "0.34%"
$cell = $ws.Cells.Item(19, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.1312"
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = "-1.56%"
$cell = $ws.Cells.Item(20, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.802"
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = "-7.38%"
$cell = $ws.Cells.Item(21, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.2618"
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = "3.97%"
$cell = $ws.Cells.Item(22, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.04416"
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = "-0.20%"
$cell = $ws.Cells.Item(23, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.001244"
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = "0.42%"
$cell = $ws.Cells.Item(24, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = "@"
$cell.Value = "-9.59%"
$cell = $ws.Cells.Item(25, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(26, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0001301"
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = "-0.36%"
$cell = $ws.Cells.Item(27, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(28, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(29, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(30, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(31, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(32, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(33, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(34, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(35, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(36, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(37, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(38, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.02100"
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = "9.45%"
$cell = $ws.Cells.Item(39, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = "1.91%"
$cell = $ws.Cells.Item(40, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.007468"
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = "-1.27%"
$cell = $ws.Cells.Item(41, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.01014"
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = "-0.54%"
$cell = $ws.Cells.Item(42, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.1366"
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = "1.61%"
$cell = $ws.Cells.Item(43, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = "-1.74%"
$cell = $ws.Cells.Item(44, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.009676"
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = "-0.53%"
$cell = $ws.Cells.Item(45, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.00006391"
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = "3.58%"
$cell = $ws.Cells.Item(46, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = "-0.45%"
$cell = $ws.Cells.Item(47, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(48, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.001601"
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = "@"
$cell.Value = "-3.97%"
$cell = $ws.Cells.Item(49, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.00002101"
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = "-0.45%"
$cell = $ws.Cells.Item(50, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0002001"
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = "-0.45%"
$cell = $ws.Cells.Item(51, 7)
$cell.NumberFormat = "@"
$cell.Value = "2"
